$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had 3 rows:
#   Row 1: a numeric "index" header row (0..10), styled (bold/border/center-top)
#   Row 2: text column headers (Issues traded, Advances, Declines, Scraped @, ...)
#   Row 3: the actual scraped data values (including a "Scraped @" date column)
#
# The commit "Got rid of row header" removes the numeric index row, turns the
# old header row into row 1 (prefixed with a "name" label column), turns the
# old data row into row 2 (prefixed with a "latestClose" label column), and
# drops the now-unneeded "Scraped @" / date column entirely.

# Drop the old data row (row 3) completely; row 2 data is rebuilt from scratch below.
$ws.Rows.Item(3).Delete()

# --- Row 1: label + headers (skipping the old "Scraped @" column) ---
$ws.Range("A1").Value2 = "name"
$ws.Range("B1").Value2 = "Issues traded"
$ws.Range("C1").Value2 = "Advances"
$ws.Range("D1").Value2 = "Declines"
$ws.Range("E1").Value2 = "Unchanged"
$ws.Range("F1").Value2 = "New highs"
$ws.Range("G1").Value2 = "New lows"
$ws.Range("H1").Value2 = "Closing Arms (TRIN)†"
$ws.Range("I1").Value2 = "Block trades"
$ws.Range("J1").Value2 = "Adv. volume"
$ws.Range("K1").Value2 = "Decl. volume"
$ws.Range("L1").Value2 = "Total volume"

# --- Row 2: label + data (skipping the old "Sep 06, 2022" scrape-date column) ---
# Force the number-like strings (e.g. "4,813", "1.26") to be stored as text,
# matching the source data, instead of being auto-parsed into numeric values.
$ws.Range("B2:L2").NumberFormat = "@"

$ws.Range("A2").Value2 = "latestClose"
$ws.Range("B2").Value2 = "4,813"
$ws.Range("C2").Value2 = "1,696"
$ws.Range("D2").Value2 = "2,811"
$ws.Range("E2").Value2 = "306"
$ws.Range("F2").Value2 = "42"
$ws.Range("G2").Value2 = "216"
$ws.Range("H2").Value2 = "1.26"
$ws.Range("I2").Value2 = "21,021"
$ws.Range("J2").Value2 = "1,359,398,333"
$ws.Range("K2").Value2 = "2,842,883,736"
$ws.Range("L2").Value2 = "4,257,416,449"

# The old numeric header row carried bold/border/centered formatting across
# A1:L1. Only the new leading "name"/"latestClose" label column should keep
# that look; the rest of the cells go back to plain/default formatting.
$ws.Range("B1:L1").ClearFormats()
$ws.Range("B2:L2").ClearFormats()

# A2 ("latestClose") should pick up the same formatting that A1 ("name") has.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
